$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has per-row "group" values repeated across every column from C
# through EO (column 145). Extend that run 21 more columns, through FJ
# (column 166), duplicating each row's existing value.
for ($r = 2; $r -le 15; $r++) {
    $srcVal = $ws.Cells.Item($r, 3).Value2
    for ($c = 146; $c -le 166; $c++) {
        if ($srcVal -eq "") {
            # Row has no data (e.g. row 10) - still materialize the cell so
            # it exists in the sheet, matching the padded-but-empty source row.
            $ws.Cells.Item($r, $c).Font.Name = "Arial"
        } else {
            $ws.Cells.Item($r, $c).Value2 = $srcVal
        }
    }
}
